$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Username) - rows 2-4
$ws.Range("A2").Value = "Person1"
$ws.Range("A3").Value = "Person2"
$ws.Range("A4").Value = "Person3"

# Column B (FirstName) - rows 2-4
$ws.Range("B2").Value = "NameName1"
$ws.Range("B3").Value = "NameName2"
$ws.Range("B4").Value = "NameName3"

# Column C (MiddleName) - rows 2-4
$ws.Range("C2").Value = "NameN1"
$ws.Range("C3").Value = "NameN2"
$ws.Range("C4").Value = "NameN3"

# Column D (LastName) - rows 2-4
$ws.Range("D2").Value = "Last1"
$ws.Range("D3").Value = "Last2"
$ws.Range("D4").Value = "Last3"

# UsertypeID updates
$ws.Range("F3").Value = 112
$ws.Range("F4").Value = 113

# Rows 5-7 - clear the username/name/contact-id columns, leave contact number & usertype cells blank but present
$ws.Range("A5:D7").ClearContents()
$ws.Range("E5:F7").ClearContents()

# Update selection to D2:D4 with active cell D2
$ws.Range("D2:D4").Select()
